# A new weekly price observation was added for Membrillo (Vega Modelo de
# Temuco). It belongs chronologically right after the existing row 108, so
# insert a fresh row at position 109 - this pushes the old rows 109-133
# down to 110-134 (row 134 is new, carrying what used to be row 133) - and
# then populate the new row 109 with the new observation's data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(109).Insert()

$ws.Range("A109").Value = 10
$ws.Range("B109").Value = "Vega Modelo de Temuco"
$ws.Range("C109").Value = "La Araucanía"
$ws.Range("D109").Value = 44642
$ws.Range("E109").Value = 9
$ws.Range("F109").Value = "Fruta"
$ws.Range("G109").Value = 100104
$ws.Range("H109").Value = "Frutos de pepita"
$ws.Range("I109").Value = 100104003
$ws.Range("J109").Value = "Membrillo"
$ws.Range("K109").Value = "Champion"
$ws.Range("L109").Value = "Primera"
$ws.Range("M109").Value = 90
$ws.Range("N109").Value = 14000
$ws.Range("O109").Value = 14000
$ws.Range("P109").Value = 14000
$ws.Range("Q109").Value = "$/bandeja 18 kilos granel"
$ws.Range("R109").Value = "Región de O'Higgins"
$ws.Range("S109").Value = 778
$ws.Range("T109").Value = 18
